# Auto update Excel log
# Appends newly-logged sensor/alert rows to four sheets (ALERTS, PIR,
# Humidity, Proximity). Plain text values (times, labels, statuses) are
# written straight to .Value. Values that Excel's smart-typing would
# otherwise reinterpret (a bare date like "2026-01-30", or a percentage
# like "86.0%") are written as Text explicitly (NumberFormat "@") and then
# have their cell formatting cleared again so the stored cell carries no
# extra style - matching a plain appended text value.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

function Add-LogRow($ws, $row, $date, $timestamp, $hour, $location, $value, $status) {
    Set-TextValue $ws.Cells.Item($row, 1) $date
    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location
    Set-TextValue $ws.Cells.Item($row, 5) $value
    $ws.Cells.Item($row, 6).Value = $status
}

# ---- ALERTS: new row 14 ----
$wsAlerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $wsAlerts 14 "2026-01-30" "17:43:57" "17:00" "Living Room" "CRITICAL" "FALL_DETECTED"

# ---- PIR: new rows 337-345 ----
$wsPir = $wb.Worksheets.Item("PIR")
$pirTimes = @("17:43:52","17:43:53","17:43:58","17:43:59","17:44:04","17:44:09","17:44:14","17:44:19","17:44:24")
$r = 337
foreach ($t in $pirTimes) {
    Add-LogRow $wsPir $r "2026-01-30" $t "17:00" "Bathroom" "No Motion" "Inactive"
    $r = $r + 1
}

# ---- Humidity: new rows 233-238 ----
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("17:43:52","86.0%"),
    @("17:43:58","87.0%"),
    @("17:44:00","86.9%"),
    @("17:44:15","86.1%"),
    @("17:44:20","86.9%"),
    @("17:44:25","86.1%")
)
$r = 233
foreach ($row in $humidityRows) {
    Add-LogRow $wsHumidity $r "2026-01-30" $row[0] "17:00" "Bathroom" $row[1] "Active"
    $r = $r + 1
}

# ---- Proximity: new row 56 ----
$wsProximity = $wb.Worksheets.Item("Proximity")
Add-LogRow $wsProximity 56 "2026-01-30" "17:43:51" "17:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
